# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) of each data table so the
# charts built from these ranges pick up a proper series/category title,
# fixes missing Portuguese accents in several labels, and refreshes a
# couple of stale figures on the "Custo Total" sheet. Also drops the
# now-unused "Teto" row from the "Emissoes Totais" sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    param($ws, [string]$addr, [string]$text)

    $ws.Range($addr).Value = $text
    # Match the existing bold/centered/bordered header style already used
    # by the rest of row 1 (B1 carries it) instead of minting a new style.
    $ws.Range("B1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

function Clear-LabelStyle {
    param($ws, [string]$addr)
    $ws.Range($addr).Style = "Normal"
}

# --- Sheets 1-4: Potencia Acumulada, Geracao Periodo Medio, Atendimento a
#     Ponta, Potencia Incremental - all share the same row layout. ---
$sheetsWithFonte = @(1, 2, 3, 4)
foreach ($idx in $sheetsWithFonte) {
    $ws = $wb.Worksheets.Item($idx)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia"

    Clear-LabelStyle $ws "A2"

    $ws.Range("A3").Value = "Gás Natural"
    Clear-LabelStyle $ws "A3"

    $ws.Range("A4").Value = "Carvão"
    Clear-LabelStyle $ws "A4"

    Clear-LabelStyle $ws "A5"

    $ws.Range("A6").Value = "Óleos Comb"
    Clear-LabelStyle $ws "A6"

    Clear-LabelStyle $ws "A7"

    $ws.Range("A8").Value = "Eólica"
    Clear-LabelStyle $ws "A8"

    Clear-LabelStyle $ws "A9"

    Clear-LabelStyle $ws "A10"

    $ws.Range("A11").Value = "Pot. Compl."
    Clear-LabelStyle $ws "A11"

    Clear-LabelStyle $ws "A12"
}

# --- Sheet 5: Emissoes Totais (MtCO2eq) ---
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período"

$ws5.Range("A2").Value = "P.Médio"
Clear-LabelStyle $ws5 "A2"

$ws5.Range("A3").Value = "P.Crítico"
Clear-LabelStyle $ws5 "A3"

# Drop the "Teto" row entirely (dimension shrinks from A1:E4 to A1:E3).
$ws5.Rows.Item(4).Delete()

# --- Sheet 6: Custo Total (bilhões de R$) ---
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "Tipo Expansão"

# B1 keeps its existing header style, but the label itself changes from
# "Custo" to the text "2015" (a year label, not a numeric value) - force
# text via a quote-prefix, then re-apply the original style so the
# quote-prefix formatting doesn't stick to the cell.
$ws6.Range("B1").Value = "'2015"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Range("A2").Value = "Expansão Centralizada"
Clear-LabelStyle $ws6 "A2"
$ws6.Range("B2").Value = 612

$ws6.Range("A3").Value = "Expansão por GD"
Clear-LabelStyle $ws6 "A3"
$ws6.Range("B3").Value = 99
